# Map sheet, row 28 ("Variable Type") currently has an empty target-column
# value (B28). The author filled it in with "VarType".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Map")

$ws.Range("B28").Value = "VarType"

# The author's saved view also leaves the cursor on E21 after making the
# edit (visible in the sheetView/<selection> of the diff).
[void]$ws.Range("E21").Select()
